$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns remain stored as text, matching the
# original inline-string cell type, so numeric-looking values such as "1.004"
# or "0.01920" are not silently reinterpreted as numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.699.29'
$ws.Range("E2").Value = '  -2.11%  '
$ws.Range("D3").Value = '1.859.81'
$ws.Range("E3").Value = '  -2.56%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '290.86'
$ws.Range("E5").Value = '  -5.50%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '0.5251'
$ws.Range("E7").Value = '  -1.03%  '
$ws.Range("D8").Value = '0.3699'
$ws.Range("E8").Value = '  -3.19%  '
$ws.Range("D9").Value = '0.07093'
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("D10").Value = '21.08'
$ws.Range("E10").Value = '  -4.57%  '
$ws.Range("D11").Value = '0.8775'
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("D12").Value = '0.08086'
$ws.Range("E12").Value = '  -1.29%  '
$ws.Range("D13").Value = '1.963.29'
$ws.Range("E13").Value = '  +64.50%  '
$ws.Range("D14").Value = '91.36'
$ws.Range("E14").Value = '  -4.58%  '
$ws.Range("D15").Value = '5.228'
$ws.Range("E15").Value = '  -2.19%  '
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").Value = '14.59'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '0.000008422'
$ws.Range("E18").Value = '  -2.72%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = '26.734.89'
$ws.Range("E20").Value = '  -2.11%  '
$ws.Range("D21").Value = '4.915'
$ws.Range("E21").Value = '  -2.75%  '
$ws.Range("D22").Value = '10.53'
$ws.Range("E22").Value = '  -2.72%  '
$ws.Range("D23").Value = '6.304'
$ws.Range("E23").Value = '  -3.29%  '
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").Value = '145.05'
$ws.Range("E24").Value = '  -3.28%  '
$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").Value = '2.235'
$ws.Range("E25").Value = '  -2.50%  '
$ws.Range("D26").Value = '1.733'
$ws.Range("E26").Value = '  -0.70%  '
$ws.Range("D27").Value = '17.80'
$ws.Range("E27").Value = '  -2.43%  '
$ws.Range("D28").Value = '112.87'
$ws.Range("E28").Value = '  -3.53%  '
$ws.Range("D29").Value = '4.649'
$ws.Range("E29").Value = '  -3.56%  '
$ws.Range("D30").Value = '4.564'
$ws.Range("E30").Value = '  -5.16%  '
$ws.Range("D31").Value = '0.09046'
$ws.Range("E31").Value = '  -2.60%  '
$ws.Range("D32").Value = '0.7902'
$ws.Range("E32").Value = '  -5.66%  '
$ws.Range("D33").Value = '0.04957'
$ws.Range("E33").Value = '  -2.19%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.950'
$ws.Range("E34").Value = '  -1.62%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.149'
$ws.Range("E35").Value = '  -6.26%  '
$ws.Range("D36").Value = '0.5856'
$ws.Range("E36").Value = '  +1.84%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '3.177'
$ws.Range("E37").Value = '  -5.36%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.573'
$ws.Range("E38").Value = '  -4.37%  '
$ws.Range("D39").Value = '0.01920'
$ws.Range("E39").Value = '  -4.26%  '
$ws.Range("D40").Value = '1.054'
$ws.Range("E40").Value = '  -2.20%  '
$ws.Range("D41").Value = '0.5147'
$ws.Range("E41").Value = '  +4.53%  '
$ws.Range("D42").Value = '6.402'
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '114.32'
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '8.502'
$ws.Range("E44").Value = '  -8.64%  '
$ws.Range("D45").Value = '0.1470'
$ws.Range("E45").Value = '  -3.54%  '
$ws.Range("D46").Value = '1.003'
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").Value = '9.958'
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("D48").Value = '1.608'
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("D49").Value = '36.76'
$ws.Range("E49").Value = '  -4.85%  '
$ws.Range("D50").Value = '0.06045'
$ws.Range("E50").Value = '  -1.72%  '
$ws.Range("D51").Value = '61.57'
$ws.Range("E51").Value = '  -3.17%  '
